# Update the StructureDefinition workbook:
#  - "Metadata" sheet: refresh URL / Version / Date / Publisher to the
#    LinuxForHealth re-branding values.
#  - "Elements" sheet: clear the stray Constraint(s) text that had been
#    duplicated onto the root "Extension" row.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/job-location-zipcode"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
